$d = $word.ActiveDocument

# 1. Shorten the first bullet: remove the DMRB clause, keep only the lead-in sentence.
$d.Content.Find.Execute(
    "Led service operations across 3 multifamily properties while designing and deploying the Make Ready Digital Board (DMRB) — a logic-based AI tool used live to coordinate unit readiness",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Led service operations across 3 multifamily properties",
    2
)

# 2. Remove the whole bullet paragraph about reducing unit turnover time (13-20 days to 7).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Reduced unit turnover time from 13*20 days to 7*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}
